# Updated burndown charts to include progress from a member
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the ACTUAL REMAINING values (column G) for rows 12 through 25
$ws.Range("G12").Value = 66
$ws.Range("G13").Value = 62
$ws.Range("G14").Value = 53
$ws.Range("G15").Value = 49
$ws.Range("G16").Value = 49
$ws.Range("G17").Value = 49
$ws.Range("G18").Value = 47
$ws.Range("G19").Value = 49
$ws.Range("G20").Value = 49
$ws.Range("G21").Value = 43
$ws.Range("G22").Value = 43
$ws.Range("G23").Value = 43
$ws.Range("G24").Value = 43
$ws.Range("G25").Value = 43

# Update the active selection to reflect the cell last edited
$ws.Range("J30").Select()
